# Auto-generated COM-interop script applying the "openApi customization" edit
# to the "Especialista Spring Rest" glossary workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Termos")

# --- Row 368: only the row height shrinks (text content is unchanged) ---
$ws.Rows.Item(368).RowHeight = 121.65

# --- Row 370 ---
$ws.Range("A370").Value = "Autorizando"
$ws.Range("B370").Value = "Ao clicar no botão ao canto da tela em documentação é possível observar que tudo está preenchido e configurado`nSelecione todos os escopos e clique em autorizar`nPreencha com as credenciais de seu usuário joao.ger@algafood.com.br"
$ws.Range("C370").NumberFormat = "@"
$ws.Range("C370").Value = "26.5"
$ws.Rows.Item(370).RowHeight = 41.55

# --- Row 371 ---
$ws.Range("A371").Value = "Criação de TAGs"
$ws.Range("B371").Value = "Servem para agrupar alguns recursos semelhantes definindo nome e descrição`n@Tag(name = `"Cidades`", description = `"Gerencia as Cidades`")`npublic interface CidadeControllerOpenApi {…}"
$ws.Range("C371").NumberFormat = "@"
$ws.Range("C371").Value = "26.6"
$ws.Rows.Item(371).RowHeight = 43.35

# --- Row 372 ---
$ws.Range("A372").Value = "Descrevendo Operações"
$ws.Range("B372").Value = "Descreve uma aoperação de um endpoint na documentação:`n`t@Operation(summary = `"Cadastra uma cidade`", description = `"Cadastro de uma cidade, necessita de um estado e um nome válido`")`n`tCidadeModel adicionar(CidadeInput cidadeInput);"
$ws.Range("C372").NumberFormat = "@"
$ws.Range("C372").Value = "26.7"
$ws.Rows.Item(372).RowHeight = 39.75

# --- Row 373 ---
$ws.Range("A373").Value = "Descrevendo Parâmetros"
$ws.Range("B373").Value = "Define os parâmetros nos recursos attravés da anotação @Parameter e também @Requestbody do swagger`nTambém é possível definir um exemplo que autocompletará na documentação:`nCidadeModel buscar(@Parameter(description = `"ID de uma cidade`", example = `"1`", required = true) Long cidadeId);`n`tCidadeModel adicionar(@RequestBody(description = `"Representação de uma nova cidade`", required = true) CidadeInput cidadeInput);"
$ws.Range("C373").NumberFormat = "@"
$ws.Range("C373").Value = "26.8"
$ws.Rows.Item(373).RowHeight = 70.45

# --- Row 374 ---
$ws.Range("A374").Value = "Descrevendo Modelos de Representação"
$ws.Range("B374").Value = "- Podemos descrever os Modelos de Representação, seus campos e classes`n- Essa configuração é refletida tanto requisição quanto nos modelos representados abaixo na documentação`n@NotBlank`n@Schema(example = `"Uberlândia`")`nprivate String nome;"
$ws.Range("C374").NumberFormat = "@"
$ws.Range("C374").Value = "26.9"
$ws.Rows.Item(374).RowHeight = 74.05

# --- Row 375 ---
$ws.Range("A375").Value = "Descrevendo Validação de Modelo"
$ws.Range("B375").Value = "- Para definir as validações de modelo o springDocs consegue aproveitar as anotações já definidas como @NotBlank, @NotNull`n- Mas podemos definir na própria anotação caso necessário com required = true`n- Porém é preferível continuar utilizando as anotações do Javax`n@Schema(example = `"Uberlândia`", required = true)`nprivate String nome;"
$ws.Range("C375").NumberFormat = "@"
$ws.Range("C375").Value = "26.10"
$ws.Rows.Item(375).RowHeight = 75.3

# --- Row 376 ---
$ws.Range("A376").Value = "Definindo Códigos de Respostas `nDe Forma Global"
$ws.Range("B376").Value = "Necessário a criação de um Bean OpenAPICustomizer`nEssa customização faz um iteração sobre os paths da nossa aplicação e adiciona ApiResponse Globeis nele de acordo com seus status`n                        ApiResponse apiResponseNaoEncontrado = new ApiResponse().description(`"Recurso não encontrado`");`n                        responses.addApiResponse(`"406`", apiResponseNaoEncontrado);"
$ws.Range("C376").NumberFormat = "@"
$ws.Range("C376").Value = "26.11"
$ws.Rows.Item(376).RowHeight = 54.2

# --- Row 377 ---
$ws.Range("A377").Value = "Descrevendo StatusCode `nPara Respostas Específicas"
$ws.Range("B377").Value = "Na anotação de @Operation basta passar o código e descrição`n`t@Operation(summary = `"Busca uma cidade por Id`", responses = {`n`t`t`t@ApiResponse(responseCode = `"200`"),`n`t`t`t@ApiResponse(responseCode = `"400`", description = `"ID da cidade inválido`",`n`t`t`t`t`tcontent = @Content(schema = @Schema))})`n`tCidadeModel buscar(@Parameter(description = `"ID de uma cidade`", example = `"1`", required = true) Long cidadeId);"
$ws.Range("C377").NumberFormat = "@"
$ws.Range("C377").Value = "26.12"
$ws.Rows.Item(377).RowHeight = 89.15

# --- Row 378 ---
$ws.Range("A378").Value = "Descrevendo StatusCode`nDe Acordo com Método HTTP"
$ws.Range("C378").NumberFormat = "@"
$ws.Range("C378").Value = "26.13"

# --- Update the view: scroll position + active selection ---
$ws.Activate()
$ws.Range("C379").Select()
try { $excel.ActiveWindow.ScrollRow = 376 } catch { }
try { $excel.ActiveWindow.ScrollColumn = 1 } catch { }

